# Apply the BOM update described in the commit:
#   - Bump the SKU code for the fixed resistor line on the Electrical sheet
#     from ...AAAA to ...AAAB.
#   - Remove the obsolete "(vide)" placeholder component line
#     (MECA-MECA-STD-VIDE-AAAA) from the Mechanical sheet; every row below
#     it shifts up by one.
#   - Refresh the remembered cursor/selection state on both sheets to
#     reflect where the user was last working.

$wb = $excel.ActiveWorkbook

$wsElec = $wb.Worksheets.Item("SKU_Électrique")
$wsMeca = $wb.Worksheets.Item("SKU_Mécanique")

# 1. Electrical sheet: update the resistor SKU suffix (AAAA -> AAAB).
$wsElec.Range("A161").Value = "ELEC-ELEC-STD-RESIST-AAAB"

# 2. Mechanical sheet: delete the empty/placeholder component row.
#    Excel automatically shifts every following row up by one and
#    recalculates the sheet dimension.
$wsMeca.Rows.Item(236).Delete()

# 3. Restore the Electrical sheet's last selected cell without making it
#    the active tab (it stays in the background).
$wsElec.Activate()
[void]$wsElec.Range("B21").Select()

# 4. Re-activate the Mechanical sheet (the tab that is actually selected
#    in the workbook) and select the full sheet, clearing out the old
#    scroll position / single-cell selection.
$wsMeca.Activate()
[void]$wsMeca.Cells.Select()
